$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on price cells whose new values would otherwise be
# auto-coerced into numbers by Excel, so they stay text like the rest of column D.
$textForceAddrs = @("D5","D6","D7","D8","D11","D12","D17","D20","D21","D22","D23","D24","D25","D27","D28","D30","D31","D33","D34","D36","D38","D39","D40","D41","D43","D44","D46","D49","D50","D51")
foreach ($addr in $textForceAddrs) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated cell values row by row
$ws.Range("D2").Value = "58.622.65"
$ws.Range("E2").Value = "  +2.00%  "
$ws.Range("D3").Value = "2.971.03"
$ws.Range("E3").Value = "  +2.36%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "557.66"
$ws.Range("E5").Value = "  +1.42%  "
$ws.Range("D6").Value = "135.99"
$ws.Range("E6").Value = "  +10.14%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").Value = "0.515"
$ws.Range("E8").Value = "  +3.56%  "
$ws.Range("D9").Value = "2.970.22"
$ws.Range("E9").Value = "  +2.47%  "
$ws.Range("E10").Value = "  +4.69%  "
$ws.Range("D11").Value = "4.81"
$ws.Range("E11").Value = "  +1.57%  "
$ws.Range("D12").Value = "0.453"
$ws.Range("E12").Value = "  +3.20%  "
$ws.Range("E13").Value = "  +6.32%  "
$ws.Range("E14").Value = "  +2.71%  "
$ws.Range("E15").Value = "  +2.69%  "
$ws.Range("D16").Value = "3.431.21"
$ws.Range("E16").Value = "  +1.88%  "
$ws.Range("D17").Value = "6.92"
$ws.Range("E17").Value = "  +4.83%  "
$ws.Range("D18").Value = "2.971.83"
$ws.Range("E18").Value = "  +2.62%  "
$ws.Range("D19").Value = "58.575.02"
$ws.Range("E19").Value = "  +1.94%  "
$ws.Range("D20").Value = "419.24"
$ws.Range("E20").Value = "  +2.93%  "
$ws.Range("D21").Value = "13.44"
$ws.Range("E21").Value = "  +3.76%  "
$ws.Range("D22").Value = "0.708"
$ws.Range("E22").Value = "  +5.31%  "
$ws.Range("D23").Value = "7.06"
$ws.Range("E23").Value = "  +3.18%  "
$ws.Range("D24").Value = "13.27"
$ws.Range("E24").Value = "  +2.94%  "
$ws.Range("D25").Value = "79.86"
$ws.Range("E25").Value = "  +3.47%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.23%  "
$ws.Range("D28").Value = "2.10"
$ws.Range("E28").Value = "  +8.32%  "
$ws.Range("E29").Value = "  +2.03%  "
$ws.Range("D30").Value = "7.70"
$ws.Range("E30").Value = "  +6.61%  "
$ws.Range("D31").Value = "25.48"
$ws.Range("E31").Value = "  +2.97%  "
$ws.Range("E32").Value = "  -0.43%  "
$ws.Range("D33").Value = "0.0983"
$ws.Range("E33").Value = "  +0.64%  "
$ws.Range("D34").Value = "0.987"
$ws.Range("E34").Value = "  +8.46%  "
$ws.Range("D35").Value = "0.0₃0751"
$ws.Range("E35").Value = "  +20.52%  "
$ws.Range("D36").Value = "5.69"
$ws.Range("E36").Value = "  +5.36%  "
$ws.Range("D38").Value = "48.58"
$ws.Range("E38").Value = "  +1.52%  "
$ws.Range("D39").Value = "8.65"
$ws.Range("E39").Value = "  +2.46%  "
$ws.Range("D40").Value = "2.73"
$ws.Range("E40").Value = "  +12.95%  "
$ws.Range("D41").Value = "396.04"
$ws.Range("E41").Value = "  +9.80%  "
$ws.Range("D42").Value = "2.737.09"
$ws.Range("E42").Value = "  +4.55%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "0.0346"
$ws.Range("E43").Value = "  +1.13%  "
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").Value = "0.107"
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").Value = "124.64"
$ws.Range("E46").Value = "  +4.10%  "
$ws.Range("E47").Value = "  +5.03%  "
$ws.Range("E48").Value = "  +2.99%  "
$ws.Range("D49").Value = "0.109"
$ws.Range("E49").Value = "  +1.83%  "
$ws.Range("B50").Value = "Arweave"
$ws.Range("C50").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D50").Value = "31.87"
$ws.Range("E50").Value = "  +19.36%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "23.08"
$ws.Range("E51").Value = "  +0.49%  "

# Reset style on the forced cells back to Normal so no stray formatting remains
foreach ($addr in $textForceAddrs) {
    $ws.Range($addr).Style = "Normal"
}
